# Rename the "_old"/"_new" column-header suffixes to the respective input
# file format versions ("_FV2304" / "_FV2310"), then turn the header row +
# data range into an Excel Table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Base (suffix-less) column names, in on-sheet order A.. for the "before"
# block (columns A-J), then "diff" (K), then the same bases again for the
# "after" block (columns L-U).
$baseNames = @(
    "Segmentname",
    "Segmentgruppe",
    "Segment",
    "Datenelement",
    "Segment ID",
    "Code",
    "Qualifier",
    "Beschreibung",
    "Bedingungsausdruck",
    "Bedingung"
)

$headers = @()
foreach ($name in $baseNames) { $headers += "$name`_FV2304" }
$headers += "diff"
foreach ($name in $baseNames) { $headers += "$name`_FV2310" }

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value2 = $headers[$i]
}

# Turn the used range into a proper table so the new headers are also
# reflected as the table's column names.
$dataRange = $ws.Range("A1:U70")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $dataRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"
$lo.TableStyle = ""

# Freeze the header row (row 1) so it stays visible while scrolling.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
